$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 99; $row++) {
    if ($row -eq 36) { continue }
    $cell = $ws.Cells.Item($row, 5)
    $cell.Value = $cell.Value2 - 1
}
